$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update player stats: coins, win, lost
$ws.Range("E2").Value = 441
$ws.Range("F2").Value = 7
$ws.Range("G2").Value = 7

# Update selection to match the last selected range A3:H4 with active cell H4
$excel.Goto($ws.Range("A3:H4"))
$ws.Range("H4").Activate()
